$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.13084733333333
$ws.Range("H2").Value = 30.392542
$ws.Range("I2").Value = 0.06120524725024945
$ws.Range("J2").Value = 0.0617054391363504
$ws.Range("M2").Value = 9.841031333333333
$ws.Range("N2").Value = 29.523094
$ws.Range("O2").Value = 0.1083017349730097
$ws.Range("P2").Value = 0.1125970533891552
$ws.Range("Q2").Value = 99.69798604054978
$ws.Range("R2").Value = 897.281874364948
$ws.Range("S2").Value = 0.006628634466654047
$ws.Range("T2").Value = 0.006947850624836912

# Row 3
$ws.Range("G3").Value = 10.13084733333333
$ws.Range("H3").Value = 30.392542
$ws.Range("I3").Value = 0.06120524725024945
$ws.Range("J3").Value = 0.0617054391363504
$ws.Range("O3").Value = 0.3504595127507141
$ws.Range("P3").Value = 0.3643589687437936
$ws.Range("Q3").Value = 322.618170601848
$ws.Range("R3").Value = 2903.563535416632
$ws.Range("S3").Value = 0.0214499611291094
$ws.Range("T3").Value = 0.02248293016960356

# Row 4
$ws.Range("G4").Value = 10.13084733333333
$ws.Range("H4").Value = 30.392542
$ws.Range("I4").Value = 0.06120524725024945
$ws.Range("J4").Value = 0.0617054391363504
$ws.Range("M4").Value = 16.16670066666667
$ws.Range("N4").Value = 48.500102
$ws.Range("O4").Value = 0.1779164877830196
$ws.Range("P4").Value = 0.1849727733236046
$ws.Range("Q4").Value = 163.7823763376982
$ws.Range("R4").Value = 1474.041387039284
$ws.Range("S4").Value = 0.0108894226246557
$ws.Range("T4").Value = 0.01141382620620163

# Row 5
$ws.Range("G5").Value = 10.13084733333333
$ws.Range("H5").Value = 30.392542
$ws.Range("I5").Value = 0.06120524725024945
$ws.Range("J5").Value = 0.0617054391363504
$ws.Range("M5").Value = 10.399077
$ws.Range("N5").Value = 20.798154
$ws.Range("O5").Value = 0.1144430947397913
$ws.Range("P5").Value = 0.07932132236322763
$ws.Range("Q5").Value = 105.351461494578
$ws.Range("R5").Value = 632.108768967468
$ws.Range("S5").Value = 0.007004517909632646
$ws.Range("T5").Value = 0.004894557029298973

# Row 6
$ws.Range("G6").Value = 10.13084733333333
$ws.Range("H6").Value = 30.392542
$ws.Range("I6").Value = 0.06120524725024945
$ws.Range("J6").Value = 0.0617054391363504
$ws.Range("M6").Value = 22.614852
$ws.Range("N6").Value = 67.84455600000001
$ws.Range("O6").Value = 0.2488791697534654
$ws.Range("P6").Value = 0.258749882180219
$ws.Range("Q6").Value = 229.107613077928
$ws.Range("R6").Value = 2061.968517701352
$ws.Range("S6").Value = 0.01523271112019766
$ws.Range("T6").Value = 0.01596627510640934

# Row 7
$ws.Range("I7").Value = 0.8155576738931257
$ws.Range("J7").Value = 0.8222227124225968
$ws.Range("M7").Value = 9.841031333333333
$ws.Range("N7").Value = 29.523094
$ws.Range("O7").Value = 0.1083017349730097
$ws.Range("P7").Value = 0.1125970533891552
$ws.Range("Q7").Value = 1328.472005914962
$ws.Range("R7").Value = 11956.24805323466
$ws.Range("S7").Value = 0.08832631105317758
$ws.Range("T7").Value = 0.09257985464842312

# Row 8
$ws.Range("I8").Value = 0.8155576738931257
$ws.Range("J8").Value = 0.8222227124225968
$ws.Range("O8").Value = 0.3504595127507141
$ws.Range("P8").Value = 0.3643589687437936
$ws.Range("S8").Value = 0.2858199450126906
$ws.Range("T8").Value = 0.2995842195760222

# Row 9
$ws.Range("I9").Value = 0.8155576738931257
$ws.Range("J9").Value = 0.8222227124225968
$ws.Range("M9").Value = 16.16670066666667
$ws.Range("N9").Value = 48.500102
$ws.Range("O9").Value = 0.1779164877830196
$ws.Range("P9").Value = 0.1849727733236046
$ws.Range("Q9").Value = 2182.394155267746
$ws.Range("R9").Value = 19641.54739740971
$ws.Range("S9").Value = 0.1451011569235542
$ws.Range("T9").Value = 0.1520888154064644

# Row 10
$ws.Range("I10").Value = 0.8155576738931257
$ws.Range("J10").Value = 0.8222227124225968
$ws.Range("M10").Value = 10.399077
$ws.Range("N10").Value = 20.798154
$ws.Range("O10").Value = 0.1144430947397913
$ws.Range("P10").Value = 0.07932132236322763
$ws.Range("Q10").Value = 1403.804358905013
$ws.Range("R10").Value = 8422.826153430078
$ws.Range("S10").Value = 0.09333494413911476
$ws.Range("T10").Value = 0.06521979282644021

# Row 11
$ws.Range("I11").Value = 0.8155576738931257
$ws.Range("J11").Value = 0.8222227124225968
$ws.Range("M11").Value = 22.614852
$ws.Range("N11").Value = 67.84455600000001
$ws.Range("O11").Value = 0.2488791697534654
$ws.Range("P11").Value = 0.258749882180219
$ws.Range("Q11").Value = 3052.850537945988
$ws.Range("R11").Value = 27475.6548415139
$ws.Range("S11").Value = 0.2029753167645886
$ws.Range("T11").Value = 0.212750029965247

# Row 12
$ws.Range("G12").Value = 4.746473666666666
$ws.Range("H12").Value = 14.239421
$ws.Range("I12").Value = 0.02867569560339488
$ws.Range("J12").Value = 0.02891004397895937
$ws.Range("M12").Value = 9.841031333333333
$ws.Range("N12").Value = 29.523094
$ws.Range("O12").Value = 0.1083017349730097
$ws.Range("P12").Value = 0.1125970533891552
$ws.Range("Q12").Value = 46.71019607650821
$ws.Range("R12").Value = 420.391764688574
$ws.Range("S12").Value = 0.003105627585405571
$ws.Range("T12").Value = 0.003255185765381712

# Row 13
$ws.Range("G13").Value = 4.746473666666666
$ws.Range("H13").Value = 14.239421
$ws.Range("I13").Value = 0.02867569560339488
$ws.Range("J13").Value = 0.02891004397895937
$ws.Range("O13").Value = 0.3504595127507141
$ws.Range("P13").Value = 0.3643589687437936
$ws.Range("Q13").Value = 151.152080449524
$ws.Range("R13").Value = 1360.368724045716
$ws.Range("S13").Value = 0.01004967030895356
$ws.Range("T13").Value = 0.01053363381051136

# Row 14
$ws.Range("G14").Value = 4.746473666666666
$ws.Range("H14").Value = 14.239421
$ws.Range("I14").Value = 0.02867569560339488
$ws.Range("J14").Value = 0.02891004397895937
$ws.Range("M14").Value = 16.16670066666667
$ws.Range("N14").Value = 48.500102
$ws.Range("O14").Value = 0.1779164877830196
$ws.Range("P14").Value = 0.1849727733236046
$ws.Range("Q14").Value = 76.73481899121576
$ws.Range("R14").Value = 690.6133709209419
$ws.Range("S14").Value = 0.005101879046490992
$ws.Range("T14").Value = 0.005347571011695492

# Row 15
$ws.Range("G15").Value = 4.746473666666666
$ws.Range("H15").Value = 14.239421
$ws.Range("I15").Value = 0.02867569560339488
$ws.Range("J15").Value = 0.02891004397895937
$ws.Range("M15").Value = 10.399077
$ws.Range("N15").Value = 20.798154
$ws.Range("O15").Value = 0.1144430947397913
$ws.Range("P15").Value = 0.07932132236322763
$ws.Range("Q15").Value = 49.35894513813899
$ws.Range("R15").Value = 296.153670828834
$ws.Range("S15").Value = 0.003281735348668735
$ws.Range("T15").Value = 0.002293182917990124

# Row 16
$ws.Range("G16").Value = 4.746473666666666
$ws.Range("H16").Value = 14.239421
$ws.Range("I16").Value = 0.02867569560339488
$ws.Range("J16").Value = 0.02891004397895937
$ws.Range("M16").Value = 22.614852
$ws.Range("N16").Value = 67.84455600000001
$ws.Range("O16").Value = 0.2488791697534654
$ws.Range("P16").Value = 0.258749882180219
$ws.Range("Q16").Value = 107.340799493564
$ws.Range("R16").Value = 966.067195442076
$ws.Range("S16").Value = 0.007136783313876015
$ws.Range("T16").Value = 0.007480470473380686

# Row 17
$ws.Range("G17").Value = 4.025238
$ws.Range("H17").Value = 8.050476
$ws.Range("I17").Value = 0.02431836932538577
$ws.Range("J17").Value = 0.01634473868084643
$ws.Range("M17").Value = 9.841031333333333
$ws.Range("N17").Value = 29.523094
$ws.Range("O17").Value = 0.1083017349730097
$ws.Range("P17").Value = 0.1125970533891552
$ws.Range("Q17").Value = 39.612493282124
$ws.Range("R17").Value = 237.674959692744
$ws.Range("S17").Value = 0.002633721589653699
$ws.Range("T17").Value = 0.001840369413879055

# Row 18
$ws.Range("G18").Value = 4.025238
$ws.Range("H18").Value = 8.050476
$ws.Range("I18").Value = 0.02431836932538577
$ws.Range("J18").Value = 0.01634473868084643
$ws.Range("O18").Value = 0.3504595127507141
$ws.Range("P18").Value = 0.3643589687437936
$ws.Range("Q18").Value = 128.184235441416
$ws.Range("R18").Value = 769.1054126484959
$ws.Range("S18").Value = 0.008522603864666609
$ws.Range("T18").Value = 0.005955352130139999

# Row 19
$ws.Range("G19").Value = 4.025238
$ws.Range("H19").Value = 8.050476
$ws.Range("I19").Value = 0.02431836932538577
$ws.Range("J19").Value = 0.01634473868084643
$ws.Range("M19").Value = 16.16670066666667
$ws.Range("N19").Value = 48.500102
$ws.Range("O19").Value = 0.1779164877830196
$ws.Range("P19").Value = 0.1849727733236046
$ws.Range("Q19").Value = 65.074817858092
$ws.Range("R19").Value = 390.448907148552
$ws.Range("S19").Value = 0.004326638858982955
$ws.Range("T19").Value = 0.003023331643045759

# Row 20
$ws.Range("G20").Value = 4.025238
$ws.Range("H20").Value = 8.050476
$ws.Range("I20").Value = 0.02431836932538577
$ws.Range("J20").Value = 0.01634473868084643
$ws.Range("M20").Value = 10.399077
$ws.Range("N20").Value = 20.798154
$ws.Range("O20").Value = 0.1144430947397913
$ws.Range("P20").Value = 0.07932132236322763
$ws.Range("Q20").Value = 41.858759905326
$ws.Range("R20").Value = 167.435039621304
$ws.Range("S20").Value = 0.002783069444622358
$ws.Range("T20").Value = 0.001296486285846136

# Row 21
$ws.Range("G21").Value = 4.025238
$ws.Range("H21").Value = 8.050476
$ws.Range("I21").Value = 0.02431836932538577
$ws.Range("J21").Value = 0.01634473868084643
$ws.Range("M21").Value = 22.614852
$ws.Range("N21").Value = 67.84455600000001
$ws.Range("O21").Value = 0.2488791697534654
$ws.Range("P21").Value = 0.258749882180219
$ws.Range("Q21").Value = 91.03016163477601
$ws.Range("R21").Value = 546.1809698086561
$ws.Range("S21").Value = 0.006052335567460152
$ws.Range("T21").Value = 0.004229199207935482

# Row 22
$ws.Range("G22").Value = 11.62680133333333
$ws.Range("H22").Value = 34.880404
$ws.Range("I22").Value = 0.0702430139278442
$ws.Range("J22").Value = 0.07081706578124704
$ws.Range("M22").Value = 9.841031333333333
$ws.Range("N22").Value = 29.523094
$ws.Range("O22").Value = 0.1083017349730097
$ws.Range("P22").Value = 0.1125970533891552
$ws.Range("Q22").Value = 114.4197162277751
$ws.Range("R22").Value = 1029.777446049976
$ws.Range("S22").Value = 0.007607440278118811
$ws.Range("T22").Value = 0.007973792936634387

# Row 23
$ws.Range("G23").Value = 11.62680133333333
$ws.Range("H23").Value = 34.880404
$ws.Range("I23").Value = 0.0702430139278442
$ws.Range("J23").Value = 0.07081706578124704
$ws.Range("O23").Value = 0.3504595127507141
$ws.Range("P23").Value = 0.3643589687437936
$ws.Range("Q23").Value = 370.2570231977759
$ws.Range("R23").Value = 3332.313208779984
$ws.Range("S23").Value = 0.0246173324352939
$ws.Range("T23").Value = 0.02580283305751657

# Row 24
$ws.Range("G24").Value = 11.62680133333333
$ws.Range("H24").Value = 34.880404
$ws.Range("I24").Value = 0.0702430139278442
$ws.Range("J24").Value = 0.07081706578124704
$ws.Range("M24").Value = 16.16670066666667
$ws.Range("N24").Value = 48.500102
$ws.Range("O24").Value = 0.1779164877830196
$ws.Range("P24").Value = 0.1849727733236046
$ws.Range("Q24").Value = 187.9670168668009
$ws.Range("R24").Value = 1691.703151801208
$ws.Range("S24").Value = 0.01249739032933576
$ws.Range("T24").Value = 0.01309922905619741

# Row 25
$ws.Range("G25").Value = 11.62680133333333
$ws.Range("H25").Value = 34.880404
$ws.Range("I25").Value = 0.0702430139278442
$ws.Range("J25").Value = 0.07081706578124704
$ws.Range("M25").Value = 10.399077
$ws.Range("N25").Value = 20.798154
$ws.Range("O25").Value = 0.1144430947397913
$ws.Range("P25").Value = 0.07932132236322763
$ws.Range("Q25").Value = 120.908002329036
$ws.Range("R25").Value = 725.448013974216
$ws.Range("S25").Value = 0.008038827897752749
$ws.Range("T25").Value = 0.005617303303652193

# Row 26
$ws.Range("G26").Value = 11.62680133333333
$ws.Range("H26").Value = 34.880404
$ws.Range("I26").Value = 0.0702430139278442
$ws.Range("J26").Value = 0.07081706578124704
$ws.Range("M26").Value = 22.614852
$ws.Range("N26").Value = 67.84455600000001
$ws.Range("O26").Value = 0.2488791697534654
$ws.Range("P26").Value = 0.258749882180219
$ws.Range("Q26").Value = 107.340799493564
$ws.Range("R26").Value = 966.067195442076
$ws.Range("S26").Value = 0.007136783313876015
$ws.Range("T26").Value = 0.007480470473380686
